$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row formatting from the last existing data row (450) down to the new rows (451:462)
$ws.Range("A450:I450").Copy() | Out-Null
$ws.Range("A451:I462").PasteSpecial(-4122) | Out-Null

# Column G uses a different (centered) style when left blank; copy that style from row 449 (G449 is blank there)
$ws.Range("G449").Copy() | Out-Null
foreach ($r in @(452,456,457,460,462)) {
    $ws.Range("G$r").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# New wellness entries for 2025-10-07 (serial date 45937)
$ws.Range("A451").Value = 45937
$ws.Range("B451").Value = "Yoann Martelat"
$ws.Range("C451").Value = 75
$ws.Range("D451").Value = 7
$ws.Range("E451").Value = 5
$ws.Range("F451").Value = 4
$ws.Range("G451").Value = "Genou"
$ws.Range("H451").Value = 5

$ws.Range("A452").Value = 45937
$ws.Range("B452").Value = "Ilan Ihaddadene"
$ws.Range("C452").Value = 75
$ws.Range("D452").Value = 8
$ws.Range("E452").Value = 6
$ws.Range("F452").Value = 0
$ws.Range("H452").Value = 10

$ws.Range("A453").Value = 45937
$ws.Range("B453").Value = "Amine Taiar"
$ws.Range("C453").Value = 75
$ws.Range("D453").Value = 4
$ws.Range("E453").Value = 8
$ws.Range("F453").Value = 7
$ws.Range("G453").Value = "Ischio"
$ws.Range("H453").Value = 5

$ws.Range("A454").Value = 45937
$ws.Range("B454").Value = "Maé Clavel"
$ws.Range("C454").Value = 75
$ws.Range("D454").Value = 6
$ws.Range("E454").Value = 6
$ws.Range("F454").Value = 4
$ws.Range("G454").Value = "Quadri ischio"
$ws.Range("H454").Value = 2

$ws.Range("A455").Value = 45937
$ws.Range("B455").Value = "Emmanuel Valey"
$ws.Range("C455").Value = 75
$ws.Range("D455").Value = 6
$ws.Range("E455").Value = 6
$ws.Range("F455").Value = 7
$ws.Range("G455").Value = "Adducteur "
$ws.Range("H455").Value = 8

$ws.Range("A456").Value = 45937
$ws.Range("B456").Value = "Mattheo Haon"
$ws.Range("C456").Value = 75
$ws.Range("D456").Value = 7
$ws.Range("E456").Value = 7
$ws.Range("F456").Value = 0
$ws.Range("H456").Value = 8

$ws.Range("A457").Value = 45937
$ws.Range("B457").Value = "Romain Thunet"
$ws.Range("C457").Value = 75
$ws.Range("D457").Value = 8
$ws.Range("E457").Value = 4
$ws.Range("F457").Value = 0
$ws.Range("H457").Value = 8

$ws.Range("A458").Value = 45937
$ws.Range("B458").Value = "Naim Dhib"
$ws.Range("C458").Value = 75
$ws.Range("D458").Value = 5
$ws.Range("E458").Value = 6
$ws.Range("F458").Value = 3
$ws.Range("G458").Value = "Genou"
$ws.Range("H458").Value = 7

$ws.Range("A459").Value = 45937
$ws.Range("B459").Value = "Yoan Zouma"
$ws.Range("C459").Value = 75
$ws.Range("D459").Value = 6
$ws.Range("E459").Value = 8
$ws.Range("F459").Value = 5
$ws.Range("G459").Value = "Cheville"
$ws.Range("H459").Value = 6

$ws.Range("A460").Value = 45937
$ws.Range("B460").Value = "Kamal Bafounta"
$ws.Range("C460").Value = 75
$ws.Range("D460").Value = 8
$ws.Range("E460").Value = 2
$ws.Range("F460").Value = 1
$ws.Range("H460").Value = 2

$ws.Range("A461").Value = 45937
$ws.Range("B461").Value = "Karim Belmahi"
$ws.Range("C461").Value = 75
$ws.Range("D461").Value = 6
$ws.Range("E461").Value = 8
$ws.Range("F461").Value = 4
$ws.Range("G461").Value = "Quadriceps "
$ws.Range("H461").Value = 10

$ws.Range("A462").Value = 45937
$ws.Range("B462").Value = "Naim Ighbane"
$ws.Range("C462").Value = 75
$ws.Range("D462").Value = 4
$ws.Range("E462").Value = 6
$ws.Range("F462").Value = 0
$ws.Range("H462").Value = 3

# Formulas: row 451 continues the existing column-I multiplication pattern (C*D)
$ws.Range("I451").Formula = "=C451*D451"

# New shared formula group for the freshly appended rows 452:462
$ws.Range("I452:I462").Formula = "=C452*D452"

$excel.Calculate()

# Update the view to reflect scrolling to the newly added rows, matching the final selection
$ws.Range("K456").Select()
$excel.ActiveWindow.ScrollRow = 429
$excel.ActiveWindow.ScrollColumn = 1

